$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("F2").Value = 1.85
$ws.Range("J2").Value = 1.26
$ws.Range("K2").Value = 980
$ws.Range("V2").Value = 1.25

# Row 3 updates (entire numeric range F3:AO3 replaced)
$ws.Range("F3").Value = 1.64
$ws.Range("G3").Value = 2.1
$ws.Range("H3").Value = 4.6
$ws.Range("I3").Value = 15.5
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.54
$ws.Range("O3").Value = 1.01
$ws.Range("P3").Value = 1.54
$ws.Range("Q3").Value = 2.06
$ws.Range("R3").Value = 1.15
$ws.Range("S3").Value = 2.98
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.04
$ws.Range("V3").Value = 1.06
$ws.Range("W3").Value = 1.92
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 4 updates
$ws.Range("F4").Value = 1.69
$ws.Range("G4").Value = 1.76
$ws.Range("J4").Value = 3.8
$ws.Range("K4").Value = 4
$ws.Range("L4").Value = 1.45
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 3.35
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 1.8
$ws.Range("Q4").Value = 2.02
$ws.Range("R4").Value = 1.29
$ws.Range("S4").Value = 3.7
$ws.Range("T4").Value = 1.99
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.17
$ws.Range("W4").Value = 2.32
$ws.Range("X4").Value = 13
$ws.Range("Y4").Value = 19
$ws.Range("Z4").Value = 48
$ws.Range("AA4").Value = 180
$ws.Range("AB4").Value = 7.6
$ws.Range("AC4").Value = 8.6
$ws.Range("AD4").Value = 24
$ws.Range("AE4").Value = 100
$ws.Range("AF4").Value = 9.6
$ws.Range("AG4").Value = 10.5
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 120
$ws.Range("AJ4").Value = 18
$ws.Range("AK4").Value = 20
$ws.Range("AL4").Value = 44
$ws.Range("AM4").Value = 180
$ws.Range("AN4").Value = 13
$ws.Range("AO4").Value = 140
